$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = "Om Patel"

$ws.Range("E7").Value = "None "
$ws.Range("F7").Value = "account_number=2004, client_number=2904, balance=1000.0, date_created=2024, 10, 4, management_fee=2.55"
$ws.Range("G7").Value = "account_number=2004, client_number=2904, balance=1000.0, date_created= 2024, 10, 4, management_fee=2.55"

$ws.Range("E8").Value = "None "
$ws.Range("F8").Value = 'account_number=2004, client_number=2904, balance=1000.0, date_created= 2024, 10, 4, management_fee="invalid"'
$ws.Range("G8").Value = "management_fee is set to 2.55"

$ws.Range("E9").Value = "account_number=2004, client_number=2904, balance=1000.0, date_created=date.today() - timedelta(days=11 * 365.25), management_fee=2.55"
$ws.Range("F9").Value = "None"
$ws.Range("G9").Value = "Service charge is 0.50"

$ws.Range("E10").Value = "account_number=2004, client_number=2904, balance=1000.0, date_created=date.today() - timedelta(days=10 * 365.25), management_fee=2.55"
$ws.Range("F10").Value = "None"
$ws.Range("G10").Value = "Service charge is 3.05"

$ws.Range("E11").Value = "account_number=2004, client_number=2904, balance=1000.0, date_created= 2024, 10, 4, management_fee=2.55"
$ws.Range("F11").Value = "None"
$ws.Range("G11").Value = "Service charge is 3.05"

$ws.Range("E12").Value = "account_number=2004, client_number=2904, balance=1000.0, date_created=date.today() - timedelta(days=11 * 365.25), management_fee=2.55"
$ws.Range("F12").Value = "None"
$ws.Range("G12").Value = "Account Number: 2004 Client Number: 2904 Balance: `$1,000.00 Date Created: {date.today() - timedelta(days=11 * 365.25)} Management Fee: Waived Account Type: Investment"

$ws.Range("E13").Value = "account_number=2004, client_number=2904, balance=1000.0, date_created= 2024, 10, 4, management_fee=2.55"
$ws.Range("F13").Value = "None"
$ws.Range("G13").Value = "Account Number: 2004 Client Number: 2904 Balance: `$1,000.00 Date Created: 2024, 10, 4, Management Fee: `$2.55 Account Type: Investment"

# Update the saved view/selection state to match the authored session.
$null = $ws.Range("I12").Select()
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
